$d = $word.ActiveDocument

# Pull the full document OOXML as a string so we can apply precise,
# surgical text-level edits (the Cell/Column width object-model setters
# in this host mis-map grid-column indices for merged-cell rows, so we
# go straight at the markup instead).
$xml = $d.Content.WordOpenXML

# 1) Resize the table's grid columns (tblGrid only - the diff does not
#    touch any individual cell's tcW).
$oldGrid = '<w:gridCol w:w="4878"/><w:gridCol w:w="3240"/><w:gridCol w:w="810"/><w:gridCol w:w="648"/>'
$newGrid = '<w:gridCol w:w="4759"/><w:gridCol w:w="3144"/><w:gridCol w:w="784"/><w:gridCol w:w="889"/>'
$xml = $xml.Replace($oldGrid, $newGrid)

# 2) Fill in the "Total number of words" compression-ratio cell.
$pattern1 = 'Total number of words(</w:t></w:r></w:p>.*?<w:tcW w:w="648" w:type="dxa"/></w:tcPr><w:p[^>]*><w:pPr><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>)</w:p>'
$replacement1 = 'Total number of words$1<w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>0.1634</w:t></w:r></w:p>'
$xml = [regex]::Replace($xml, $pattern1, $replacement1)

# 3) Fill in the "Total number of sentences" compression-ratio cell.
$pattern2 = 'Total number of sentences(</w:t></w:r></w:p>.*?<w:tcW w:w="648" w:type="dxa"/></w:tcPr><w:p[^>]*><w:pPr><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr></w:pPr>)</w:p>'
$replacement2 = 'Total number of sentences$1<w:r><w:rPr><w:b/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>0.3750</w:t></w:r></w:p>'
$xml = [regex]::Replace($xml, $pattern2, $replacement2)

# Write the modified package back over the whole document.
$d.Content.InsertXML($xml)
